$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2024-03-29 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-03-30 Saturday", 2) | Out-Null

# Update the division-problem table cells (Table 1, row/col are 1-indexed)
$tbl = $d.Tables.Item(1)

$tbl.Cell(1, 1).Range.Text = "43÷4=10, 3"
$tbl.Cell(1, 2).Range.Text = "93÷7=13, 2"
$tbl.Cell(1, 3).Range.Text = "97÷3=32, 1"
$tbl.Cell(1, 4).Range.Text = "36÷5=7, 1"
$tbl.Cell(1, 5).Range.Text = "64÷7=9, 1"
$tbl.Cell(5, 1).Range.Text = "27÷3=9, 0"
$tbl.Cell(5, 2).Range.Text = "64÷3=21, 1"
$tbl.Cell(5, 3).Range.Text = "44÷6=7, 2"
$tbl.Cell(5, 4).Range.Text = "27÷4=6, 3"
$tbl.Cell(5, 5).Range.Text = "59÷9=6, 5"
$tbl.Cell(9, 1).Range.Text = "41÷3=13, 2"
$tbl.Cell(9, 2).Range.Text = "35÷5=7, 0"
$tbl.Cell(9, 3).Range.Text = "21÷9=2, 3"
$tbl.Cell(9, 4).Range.Text = "58÷3=19, 1"
$tbl.Cell(9, 5).Range.Text = "62÷2=31, 0"
$tbl.Cell(13, 1).Range.Text = "92÷2=46, 0"
$tbl.Cell(13, 2).Range.Text = "78÷2=39, 0"
$tbl.Cell(13, 3).Range.Text = "58÷7=8, 2"
$tbl.Cell(13, 4).Range.Text = "19÷3=6, 1"
$tbl.Cell(13, 5).Range.Text = "42÷5=8, 2"
$tbl.Cell(17, 1).Range.Text = "13÷3=4, 1"
$tbl.Cell(17, 2).Range.Text = "37÷9=4, 1"
$tbl.Cell(17, 3).Range.Text = "19÷4=4, 3"
$tbl.Cell(17, 4).Range.Text = "76÷7=10, 6"
$tbl.Cell(17, 5).Range.Text = "46÷9=5, 1"
